$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "40÷4=10, 0"
$t.Cell(1, 2).Range.Text = "92÷7=13, 1"
$t.Cell(1, 3).Range.Text = "48÷4=12, 0"
$t.Cell(1, 4).Range.Text = "82÷9=9, 1"
$t.Cell(1, 5).Range.Text = "57÷9=6, 3"

$t.Cell(5, 1).Range.Text = "11÷7=1, 4"
$t.Cell(5, 2).Range.Text = "22÷3=7, 1"
$t.Cell(5, 3).Range.Text = "53÷6=8, 5"
$t.Cell(5, 4).Range.Text = "83÷7=11, 6"
$t.Cell(5, 5).Range.Text = "69÷9=7, 6"

$t.Cell(9, 1).Range.Text = "52÷2=26, 0"
$t.Cell(9, 2).Range.Text = "37÷7=5, 2"
$t.Cell(9, 3).Range.Text = "67÷6=11, 1"
$t.Cell(9, 4).Range.Text = "41÷5=8, 1"
$t.Cell(9, 5).Range.Text = "22÷8=2, 6"

$t.Cell(13, 1).Range.Text = "78÷8=9, 6"
$t.Cell(13, 2).Range.Text = "27÷4=6, 3"
$t.Cell(13, 3).Range.Text = "60÷6=10, 0"
$t.Cell(13, 4).Range.Text = "20÷8=2, 4"
$t.Cell(13, 5).Range.Text = "12÷6=2, 0"

$t.Cell(17, 1).Range.Text = "68÷8=8, 4"
$t.Cell(17, 2).Range.Text = "30÷6=5, 0"
$t.Cell(17, 3).Range.Text = "51÷6=8, 3"
$t.Cell(17, 4).Range.Text = "73÷4=18, 1"
$t.Cell(17, 5).Range.Text = "72÷3=24, 0"
